# Auto-generated edit script: apply updated price/profit figures
# from the refreshed market-board snapshot across all class sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 146.36
$ws.Range("I15").Value = 146.36
$ws.Range("K15").Value = 439.08
$ws.Range("M15").Value = -270.08
$ws.Range("H33").Value = 423.05405
$ws.Range("I33").Value = 153.64516
$ws.Range("J33").Value = 1815
$ws.Range("K33").Value = 153.64516
$ws.Range("L33").Value = 1815
$ws.Range("M33").Value = 75.35484
$ws.Range("N33").Value = -2273
$ws.Range("H51").Value = 26625
$ws.Range("I51").Value = 800
$ws.Range("J51").Value = 35233.332
$ws.Range("K51").Value = 800
$ws.Range("L51").Value = 35233.332
$ws.Range("M51").Value = -316
$ws.Range("N51").Value = -36201.332
$ws.Range("H92").Value = 17861764
$ws.Range("I92").Value = 23815428
$ws.Range("J92").Value = 772.8570999999999
$ws.Range("K92").Value = 23815428
$ws.Range("L92").Value = 772.8570999999999
$ws.Range("M92").Value = -23814180
$ws.Range("N92").Value = -3268.8571
$ws.Range("H96").Value = 1166.6666
$ws.Range("I96").Value = 1000
$ws.Range("J96").Value = 1500
$ws.Range("K96").Value = 3000
$ws.Range("L96").Value = 4500
$ws.Range("M96").Value = -1627
$ws.Range("N96").Value = -7246
$ws.Range("H98").Value = 4001305
$ws.Range("I98").Value = 1028.5
$ws.Range("J98").Value = 33336666
$ws.Range("K98").Value = 1028.5
$ws.Range("L98").Value = 33336666
$ws.Range("M98").Value = 469.5
$ws.Range("N98").Value = -33339662
$ws.Range("H100").Value = 6526.8184
$ws.Range("I100").Value = 2165.8333
$ws.Range("K100").Value = 2165.8333
$ws.Range("M100").Value = -1624.8333
$ws.Range("H122").Value = 4001305
$ws.Range("I122").Value = 1028.5
$ws.Range("J122").Value = 33336666
$ws.Range("K122").Value = 3085.5
$ws.Range("L122").Value = 100009998
$ws.Range("M122").Value = -635.5
$ws.Range("N122").Value = -100014898
$ws.Range("H135").Value = 9091897
$ws.Range("I135").Value = 315.9375
$ws.Range("J135").Value = 21741054
$ws.Range("K135").Value = 2843.4375
$ws.Range("L135").Value = 195669486
$ws.Range("M135").Value = -308.4375
$ws.Range("N135").Value = -195674556
$ws.Range("H137").Value = 29667.115
$ws.Range("I137").Value = 44492.434
$ws.Range("J137").Value = 1251.9166
$ws.Range("K137").Value = 133477.302
$ws.Range("L137").Value = 3755.7498
$ws.Range("M137").Value = -130927.302
$ws.Range("N137").Value = -8855.7498
$ws.Range("H138").Value = 1194.55
$ws.Range("I138").Value = 587.4912
$ws.Range("J138").Value = 1999.2559
$ws.Range("K138").Value = 1762.4736
$ws.Range("L138").Value = 5997.7677
$ws.Range("M138").Value = 3377.5264
$ws.Range("N138").Value = -16277.7677

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 48515.2
$ws.Range("J76").Value = 48515.2
$ws.Range("L76").Value = 48515.2
$ws.Range("N76").Value = -49191.2
$ws.Range("H79").Value = 48515.2
$ws.Range("J79").Value = 48515.2
$ws.Range("L79").Value = 48515.2
$ws.Range("N79").Value = -50855.2
$ws.Range("H88").Value = 4804216.5
$ws.Range("I88").Value = 1600.7142
$ws.Range("J88").Value = 13208794
$ws.Range("K88").Value = 1600.7142
$ws.Range("L88").Value = 13208794
$ws.Range("M88").Value = -1194.7142
$ws.Range("N88").Value = -13209606
$ws.Range("H91").Value = 4804216.5
$ws.Range("I91").Value = 1600.7142
$ws.Range("J91").Value = 13208794
$ws.Range("K91").Value = 1600.7142
$ws.Range("L91").Value = 13208794
$ws.Range("M91").Value = -196.7141999999999
$ws.Range("N91").Value = -13211602
$ws.Range("H132").Value = 5611406
$ws.Range("I132").Value = 6808187
$ws.Range("K132").Value = 20424561
$ws.Range("M132").Value = -20422031

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 335162.1
$ws.Range("I86").Value = 1888.8889
$ws.Range("J86").Value = 585117
$ws.Range("K86").Value = 1888.8889
$ws.Range("L86").Value = 585117
$ws.Range("M86").Value = -765.8888999999999
$ws.Range("N86").Value = -587363
$ws.Range("H89").Value = 335162.1
$ws.Range("I89").Value = 1888.8889
$ws.Range("J89").Value = 585117
$ws.Range("K89").Value = 9444.4445
$ws.Range("L89").Value = 2925585
$ws.Range("M89").Value = -3828.4445
$ws.Range("N89").Value = -2936817
$ws.Range("H132").Value = 333393600
$ws.Range("J132").Value = 333393600
$ws.Range("L132").Value = 333393600
$ws.Range("N132").Value = -333403720
$ws.Range("H134").Value = 24245.834
$ws.Range("I134").Value = 1092
$ws.Range("J134").Value = 93707.336
$ws.Range("K134").Value = 3276
$ws.Range("L134").Value = 281122.008
$ws.Range("M134").Value = -741
$ws.Range("N134").Value = -286192.008

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 22410.75
$ws.Range("J28").Value = 22410.75
$ws.Range("L28").Value = 22410.75
$ws.Range("N28").Value = -22900.75
$ws.Range("H31").Value = 8674.787
$ws.Range("I31").Value = 6338.184
$ws.Range("J31").Value = 18540.445
$ws.Range("K31").Value = 6338.184
$ws.Range("L31").Value = 18540.445
$ws.Range("M31").Value = -6043.184
$ws.Range("N31").Value = -19130.445
$ws.Range("H34").Value = 8674.787
$ws.Range("I34").Value = 6338.184
$ws.Range("J34").Value = 18540.445
$ws.Range("K34").Value = 6338.184
$ws.Range("L34").Value = 18540.445
$ws.Range("M34").Value = -6136.184
$ws.Range("N34").Value = -18944.445
$ws.Range("H58").Value = 2218.818
$ws.Range("I58").Value = 576.5
$ws.Range("J58").Value = 6598.3335
$ws.Range("K58").Value = 576.5
$ws.Range("L58").Value = 6598.3335
$ws.Range("M58").Value = -373.5
$ws.Range("N58").Value = -7004.3335
$ws.Range("H62").Value = 2608.3333
$ws.Range("I62").Value = 2633.3333
$ws.Range("K62").Value = 2633.3333
$ws.Range("M62").Value = -2009.3333
$ws.Range("H63").Value = 60135.5
$ws.Range("J63").Value = 60135.5
$ws.Range("L63").Value = 60135.5
$ws.Range("N63").Value = -61507.5
$ws.Range("H65").Value = 2608.3333
$ws.Range("I65").Value = 2633.3333
$ws.Range("K65").Value = 13166.6665
$ws.Range("M65").Value = -10046.6665
$ws.Range("H66").Value = 60135.5
$ws.Range("J66").Value = 60135.5
$ws.Range("L66").Value = 180406.5
$ws.Range("N66").Value = -187270.5
$ws.Range("H74").Value = 13089.546
$ws.Range("J74").Value = 14370
$ws.Range("L74").Value = 14370
$ws.Range("N74").Value = -16118
$ws.Range("H77").Value = 13089.546
$ws.Range("J77").Value = 14370
$ws.Range("L77").Value = 43110
$ws.Range("N77").Value = -51846
$ws.Range("H99").Value = 3668.6956
$ws.Range("I99").Value = 3282.353
$ws.Range("J99").Value = 4763.3335
$ws.Range("K99").Value = 3282.353
$ws.Range("L99").Value = 4763.3335
$ws.Range("M99").Value = -1784.353
$ws.Range("N99").Value = -7759.3335
$ws.Range("H126").Value = 3668.6956
$ws.Range("I126").Value = 3282.353
$ws.Range("J126").Value = 4763.3335
$ws.Range("K126").Value = 9847.059000000001
$ws.Range("L126").Value = 14290.0005
$ws.Range("M126").Value = -7377.059000000001
$ws.Range("N126").Value = -19230.0005
$ws.Range("H136").Value = 2218.818
$ws.Range("I136").Value = 576.5
$ws.Range("J136").Value = 6598.3335
$ws.Range("K136").Value = 1729.5
$ws.Range("L136").Value = 19795.0005
$ws.Range("M136").Value = 820.5
$ws.Range("N136").Value = -24895.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3797.0938
$ws.Range("J5").Value = 8990.416999999999
$ws.Range("L5").Value = 26971.251
$ws.Range("N5").Value = -27195.251
$ws.Range("H92").Value = 775
$ws.Range("I92").Value = 875
$ws.Range("J92").Value = 675
$ws.Range("K92").Value = 2625
$ws.Range("L92").Value = 2025
$ws.Range("M92").Value = -1377
$ws.Range("N92").Value = -4521
$ws.Range("H122").Value = 436.84616
$ws.Range("I122").Value = 283.1111
$ws.Range("J122").Value = 518.2353000000001
$ws.Range("K122").Value = 2547.9999
$ws.Range("L122").Value = 4664.117700000001
$ws.Range("M122").Value = -97.99990000000025
$ws.Range("N122").Value = -9564.117700000001
$ws.Range("H135").Value = 3797.0938
$ws.Range("J135").Value = 8990.416999999999
$ws.Range("L135").Value = 80913.753
$ws.Range("N135").Value = -85983.753

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3337327.5
$ws.Range("I70").Value = 5266846.5
$ws.Range("J70").Value = 4521.636
$ws.Range("K70").Value = 5266846.5
$ws.Range("L70").Value = 4521.636
$ws.Range("M70").Value = -5266576.5
$ws.Range("N70").Value = -5061.636
$ws.Range("H73").Value = 3337327.5
$ws.Range("I73").Value = 5266846.5
$ws.Range("J73").Value = 4521.636
$ws.Range("K73").Value = 5266846.5
$ws.Range("L73").Value = 4521.636
$ws.Range("M73").Value = -5265910.5
$ws.Range("N73").Value = -6393.636

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 264867.97
$ws.Range("I136").Value = 436162.12
$ws.Range("K136").Value = 1308486.36
$ws.Range("M136").Value = -1305936.36

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 15851.286
$ws.Range("I132").Value = 1363.5
$ws.Range("J132").Value = 35168.332
$ws.Range("K132").Value = 4090.5
$ws.Range("L132").Value = 105504.996
$ws.Range("M132").Value = -1560.5
$ws.Range("N132").Value = -110564.996
$ws.Range("H136").Value = 1940305.6
$ws.Range("I136").Value = 2041930.2
$ws.Range("J136").Value = 1432182.8
$ws.Range("K136").Value = 6125790.6
$ws.Range("L136").Value = 4296548.4
$ws.Range("M136").Value = -6123240.6
$ws.Range("N136").Value = -4301648.4

